$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 37; existing rows 37-87 shift down to 38-88.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new record's data.
$ws.Range("A37").Value = 11
$ws.Range("B37").Value = "Vega Monumental Concepción"
$ws.Range("C37").Value = "Bíobío"
$ws.Range("D37").Value = 44638
$ws.Range("E37").Value = 8
$ws.Range("F37").Value = 100112021
$ws.Range("G37").Value = "Ají"
$ws.Range("H37").Value = "Americana (o)"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 220
$ws.Range("K37").Value = 25000
$ws.Range("L37").Value = 26000
$ws.Range("M37").Value = 25455
$ws.Range("N37").Value = "`$/caja 25 kilos"
$ws.Range("O37").Value = "Región Metropolitana"
$ws.Range("P37").Value = 1018
$ws.Range("Q37").Value = 25
$ws.Range("R37").Value = "Hortaliza"
